$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$notes = @{
    3 = @{ mode = 'append'; note = 'Beschrijving in CMA/1/A.1 - Type Boringen en uitvoeringsmethodes gehanteerd bij staalname (https://reflabos.vito.be/2018/CMA_1_A.1.pdf).' }
    4 = @{ mode = 'append'; note = 'Beschrijving Techniek Verticale boringen (https://smetgroup.be/wp-content/uploads/techniek_verticale_boringen.pdf).' }
    5 = @{ mode = 'append'; note = 'Beschrijving: Boren voor de installatie van verticale bodemwarmtewisselaars (versie 3.1) - Smart Geotherm: (https://www.smartgeotherm.be/documents/2014/01/boren-voor-de-installatie-van-verticale-bodemwarmtewisselaars-bij-ondiepe-geothermie.pdf/). Boringen en exploitatie van grondwaterwinningsputten. Deel 1 Overzicht van de meest gebruikte boortechnieken (https://navigator.emis.vito.be/light?woId=76828&woLang=nl&woVersion=2024-04-07).' }
    7 = @{ mode = 'prepend'; note = 'Beschrijving in CMA/1/A.1 - Type Boringen en uitvoeringsmethodes gehanteerd bij staalname (https://reflabos.vito.be/2018/CMA_1_A.1.pdf).' }
    9 = @{ mode = 'prepend'; note = 'Beschrijving in CMA/1/A.1 - Type Boringen en uitvoeringsmethodes gehanteerd bij staalname (https://reflabos.vito.be/2018/CMA_1_A.1.pdf).' }
    11 = @{ mode = 'append'; note = 'Boringen en exploitatie van grondwaterwinningsputten. Deel 1 Overzicht van de meest gebruikte boortechnieken (https://navigator.emis.vito.be/light?woId=76828&woLang=nl&woVersion=2024-04-07).' }
    12 = @{ mode = 'prepend'; note = 'Beschrijving in CMA/1/A.1 - Type Boringen en uitvoeringsmethodes gehanteerd bij staalname (https://reflabos.vito.be/2018/CMA_1_A.1.pdf).' }
    13 = @{ mode = 'prepend'; note = 'Beschrijving in CMA/1/A.1 - Type Boringen en uitvoeringsmethodes gehanteerd bij staalname (https://reflabos.vito.be/2018/CMA_1_A.1.pdf).' }
    14 = @{ mode = 'prepend'; note = 'Boringen en exploitatie van grondwaterwinningsputten. Deel 1 Overzicht van de meest gebruikte boortechnieken (https://navigator.emis.vito.be/light?woId=76828&woLang=nl&woVersion=2024-04-07).' }
    15 = @{ mode = 'append'; note = 'Definitie overgenomen van Basis Registratie Ondergrond (https://docs.geostandaarden.nl/bro/def-im-bhr-g-20230607/#detail_class_Model_Boortechniek).' }
    16 = @{ mode = 'prepend'; note = 'Definitie overgenomen van Basis Registratie Ondergrond (https://docs.geostandaarden.nl/bro/def-im-bhr-g-20230607/#detail_class_Model_Boortechniek).' }
    17 = @{ mode = 'prepend'; note = 'Beschrijving in CMA/1/A.1 - Type Boringen en uitvoeringsmethodes gehanteerd bij staalname (https://reflabos.vito.be/2018/CMA_1_A.1.pdf).' }
    18 = @{ mode = 'prepend'; note = 'Beschrijving in CMA/1/A.1 - Type Boringen en uitvoeringsmethodes gehanteerd bij staalname (https://reflabos.vito.be/2018/CMA_1_A.1.pdf).' }
    19 = @{ mode = 'prepend'; note = 'Beschrijving in CMA/1/A.1 - Type Boringen en uitvoeringsmethodes gehanteerd bij staalname (https://reflabos.vito.be/2018/CMA_1_A.1.pdf).' }
    22 = @{ mode = 'prepend'; note = 'Beschrijving: Boren voor de installatie van verticale bodemwarmtewisselaars (versie 3.1) - Smart Geotherm: (https://www.smartgeotherm.be/documents/2014/01/boren-voor-de-installatie-van-verticale-bodemwarmtewisselaars-bij-ondiepe-geothermie.pdf/). Boringen en exploitatie van grondwaterwinningsputten. Deel 1 Overzicht van de meest gebruikte boortechnieken (https://navigator.emis.vito.be/light?woId=76828&woLang=nl&woVersion=2024-04-07).' }
    25 = @{ mode = 'prepend'; note = 'Beschrijving: Boren voor de installatie van verticale bodemwarmtewisselaars (versie 3.1) - Smart Geotherm: (https://www.smartgeotherm.be/documents/2014/01/boren-voor-de-installatie-van-verticale-bodemwarmtewisselaars-bij-ondiepe-geothermie.pdf/). Boringen en exploitatie van grondwaterwinningsputten. Deel 1 Overzicht van de meest gebruikte boortechnieken (https://navigator.emis.vito.be/light?woId=76828&woLang=nl&woVersion=2024-04-07).' }
    27 = @{ mode = 'append'; note = 'Beschrijving: Boren voor de installatie van verticale bodemwarmtewisselaars (versie 3.1) - Smart Geotherm: (https://www.smartgeotherm.be/documents/2014/01/boren-voor-de-installatie-van-verticale-bodemwarmtewisselaars-bij-ondiepe-geothermie.pdf/).' }
    28 = @{ mode = 'prepend'; note = 'Beschrijving in CMA/1/A.1 - Type Boringen en uitvoeringsmethodes gehanteerd bij staalname (https://reflabos.vito.be/2018/CMA_1_A.1.pdf).' }
    29 = @{ mode = 'append'; note = 'Bron: Afdeling Geotechniek (VO). Beschrijving in CMA/1/A.1 - Type Boringen en uitvoeringsmethodes gehanteerd bij staalname (https://reflabos.vito.be/2018/CMA_1_A.1.pdf).' }
    30 = @{ mode = 'prepend'; note = 'Beschrijving in CMA/1/A.1 - Type Boringen en uitvoeringsmethodes gehanteerd bij staalname (https://reflabos.vito.be/2018/CMA_1_A.1.pdf).' }
    31 = @{ mode = 'prepend'; note = 'Beschrijving Techniek Verticale boringen (https://smetgroup.be/wp-content/uploads/techniek_verticale_boringen.pdf).' }
    32 = @{ mode = 'append'; note = 'Beschrijving: Boren voor de installatie van verticale bodemwarmtewisselaars (versie 3.1) - Smart Geotherm: (https://www.smartgeotherm.be/documents/2014/01/boren-voor-de-installatie-van-verticale-bodemwarmtewisselaars-bij-ondiepe-geothermie.pdf).' }
    33 = @{ mode = 'append'; note = 'Beschrijving in CMA/1/A.1 - Type Boringen en uitvoeringsmethodes gehanteerd bij staalname (https://reflabos.vito.be/2018/CMA_1_A.1.pdf). Beschrijving: Boren voor de installatie van verticale bodemwarmtewisselaars (versie 3.1) - Smart Geotherm: (https://www.smartgeotherm.be/documents/2014/01/boren-voor-de-installatie-van-verticale-bodemwarmtewisselaars-bij-ondiepe-geothermie.pdf).' }
    35 = @{ mode = 'prepend'; note = 'Beschrijving in CMA/1/A.1 - Type Boringen en uitvoeringsmethodes gehanteerd bij staalname (https://reflabos.vito.be/2018/CMA_1_A.1.pdf).' }
    36 = @{ mode = 'prepend'; note = 'Beschrijving: Boren voor de installatie van verticale bodemwarmtewisselaars (versie 3.1) - Smart Geotherm: (https://www.smartgeotherm.be/documents/2014/01/boren-voor-de-installatie-van-verticale-bodemwarmtewisselaars-bij-ondiepe-geothermie.pdf/). Boringen en exploitatie van grondwaterwinningsputten. Deel 1 Overzicht van de meest gebruikte boortechnieken (https://navigator.emis.vito.be/light?woId=76828&woLang=nl&woVersion=2024-04-07).' }
    37 = @{ mode = 'prepend'; note = 'Beschrijving in CMA/1/A.1 - Type Boringen en uitvoeringsmethodes gehanteerd bij staalname (https://reflabos.vito.be/2018/CMA_1_A.1.pdf).' }
    38 = @{ mode = 'append'; note = 'Beschrijving in CMA/1/A.4 - Apparatuur en materiaal (https://reflabos.vito.be/2024/CMA_1_A.4.pdf).' }
    39 = @{ mode = 'prepend'; note = 'Beschrijving in CMA/1/A.4 - Apparatuur en materiaal (https://reflabos.vito.be/2024/CMA_1_A.4.pdf).' }
    40 = @{ mode = 'append'; note = 'Beschrijving in CMA/1/A.1 - Type Boringen en uitvoeringsmethodes gehanteerd bij staalname (https://reflabos.vito.be/2018/CMA_1_A.1.pdf). Beschrijving in CMA/1/A.4 - Apparatuur en materiaal (https://reflabos.vito.be/2024/CMA_1_A.4.pdf).' }
}

foreach ($row in $notes.Keys) {
    $entry = $notes[$row]
    $cell = $ws.Cells.Item([int]$row, 12)
    $current = $cell.Value2
    if ($entry.mode -eq "append") {
        $cell.Value2 = $current + "|" + $entry.note
    } else {
        $cell.Value2 = $entry.note + "|" + $current
    }
}
